$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# StatQuery text (column C on the SamplesTab/FilesTab rows) picks up an extra
# space before the bracket: "in  [" / previously "in [".
$statQueryNew = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen in  ["Taxane only"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

# SamplesTab Cypher query (column B, row 3) gets the same extra-space fix:
# "IN  [" instead of "IN [".
$samplesQueryNew = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN  ["Taxane only"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
 order By samp.sample_id ASC LIMIT 100
'@

# Apply the StatQuery fix first (rows 3 and 4, column C) ...
$ws.Range("C3").Value = $statQueryNew
$ws.Range("C4").Value = $statQueryNew

# ... then the SamplesTab query fix (row 3, column B). Column B row 4 (FilesTab
# query) keeps its original text -- only its shared-string slot shifts because
# the old, now-unreferenced SamplesTab string is dropped from the pool.
$ws.Range("B3").Value = $samplesQueryNew

# New row 5: a single formatted (wrap-text styled) but otherwise empty cell in
# C5, matching the rest of column C's wrap-text style.
$ws.Range("C5").WrapText = $true

# Row heights grew slightly to accommodate the longer wrapped text.
$ws.Rows(2).RowHeight = 345
$ws.Rows(3).RowHeight = 375
$ws.Rows(4).RowHeight = 409.5

# Selection / scroll position moved up one row (now centered on the edited
# SamplesTab row).
$ws.Range("B3").Select() | Out-Null
